$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This edit:
#   1) Swaps the full contents of rows 2 and 3 (case A 13467-2023 / A 45325-2025)
#   2) Swaps the full contents of rows 6 and 7 (case A 12651-2022 / A 5792-2024)
#   3) Rotates rows 12 -> 13 -> 14 -> 12 (cases A 13651-2023 / A 8194-2025 / A 50997-2025)
#   4) Updates column C ("Förändrad") for every data row (2..16) from 46079 to 46081
# ---------------------------------------------------------------------------

# ---- Row 2 becomes old row 3 (A 45325-2025) -------------------------------
$ws.Range("A2").Value = "A 45325-2025"
$ws.Range("B2").Value = 45922
$ws.Range("C2").Value = 46081
$ws.Range("D2").Value = "SKÅNE LÄN"
$ws.Range("E2").Value = "HELSINGBORG"
$ws.Range("G2").Value = 1.6
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = "Nordlig buksimmare`r`nStörre vattensalamander`r`nÅkergroda`r`nMindre vattensalamander`r`nVanlig groda"
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/artfynd/A 45325-2025 artfynd.xlsx", "A 45325-2025")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/kartor/A 45325-2025 karta.png", "A 45325-2025")'
$ws.Range("U2").Value = ""
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/klagomål/A 45325-2025 FSC-klagomål.docx", "A 45325-2025")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/klagomålsmail/A 45325-2025 FSC-klagomål mail.docx", "A 45325-2025")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/tillsyn/A 45325-2025 tillsynsbegäran.docx", "A 45325-2025")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/tillsynsmail/A 45325-2025 tillsynsbegäran mail.docx", "A 45325-2025")'
$ws.Range("Z2").ClearContents()

# ---- Row 3 becomes old row 2 (A 13467-2023) -------------------------------
$ws.Range("A3").Value = "A 13467-2023"
$ws.Range("B3").Value = 45005
$ws.Range("C3").Value = 46081
$ws.Range("D3").Value = "SKÅNE LÄN"
$ws.Range("E3").Value = "HELSINGBORG"
$ws.Range("G3").Value = 2.3
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 5
$ws.Range("R3").Value = "Gulsparv`r`nHypoxylon petriniae`r`nKråka`r`nGrå skärelav`r`nGulnål"
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/artfynd/A 13467-2023 artfynd.xlsx", "A 13467-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/kartor/A 13467-2023 karta.png", "A 13467-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/klagomål/A 13467-2023 FSC-klagomål.docx", "A 13467-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/klagomålsmail/A 13467-2023 FSC-klagomål mail.docx", "A 13467-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/tillsyn/A 13467-2023 tillsynsbegäran.docx", "A 13467-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/tillsynsmail/A 13467-2023 tillsynsbegäran mail.docx", "A 13467-2023")'
$ws.Range("Z3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1283/fåglar/A 13467-2023 prioriterade fågelarter.docx", "A 13467-2023")'

# ---- Row 4 unchanged except "Förändrad" date -------------------------------
$ws.Range("C4").Value = 46081

# ---- Row 5 unchanged except "Förändrad" date -------------------------------
$ws.Range("C5").Value = 46081

# ---- Row 6 becomes old row 7 (A 5792-2024) --------------------------------
$ws.Range("A6").Value = "A 5792-2024"
$ws.Range("B6").Value = 45335
$ws.Range("C6").Value = 46081
$ws.Range("G6").Value = 5.6

# ---- Row 7 becomes old row 6 (A 12651-2022) -------------------------------
$ws.Range("A7").Value = "A 12651-2022"
$ws.Range("B7").Value = 44641
$ws.Range("C7").Value = 46081
$ws.Range("G7").Value = 3.2

# ---- Rows 8..11 unchanged except "Förändrad" date --------------------------
$ws.Range("C8").Value = 46081
$ws.Range("C9").Value = 46081
$ws.Range("C10").Value = 46081
$ws.Range("C11").Value = 46081

# ---- Row 12 becomes old row 13 (A 8194-2025) -------------------------------
$ws.Range("A12").Value = "A 8194-2025"
$ws.Range("B12").Value = 45708
$ws.Range("C12").Value = 46081
$ws.Range("G12").Value = 1.9

# ---- Row 13 becomes old row 14 (A 50997-2025) ------------------------------
$ws.Range("A13").Value = "A 50997-2025"
$ws.Range("B13").Value = 45946
$ws.Range("C13").Value = 46081
$ws.Range("G13").Value = 1.5

# ---- Row 14 becomes old row 12 (A 13651-2023) ------------------------------
$ws.Range("A14").Value = "A 13651-2023"
$ws.Range("B14").Value = 45006
$ws.Range("C14").Value = 46081
$ws.Range("G14").Value = 2.2

# ---- Rows 15..16 unchanged except "Förändrad" date --------------------------
$ws.Range("C15").Value = 46081
$ws.Range("C16").Value = 46081
